# "Lots of updates for Sirenum, Out There"
#
# Update the input numbers on Sheet1 (the "ROTATIONAL GRAVITY IN A HABITAT"
# and "Thrust" tables). All the other changed cells in this sheet are
# formulas that depend on these inputs, so they recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# ROTATIONAL GRAVITY IN A HABITAT: Radius (G4) and rpm (G5)
$ws.Range("G4").Value = 40
$ws.Range("G5").Value = 1.2

# Thrust section: G (B20)
$ws.Range("B20").Value = 0.0013888888888888889

# Recalculate every dependent formula (G8:G10, B22:B35, etc.)
$excel.CalculateFullRebuild()

# Restore the view: scroll position and current selection
$window = $excel.ActiveWindow
$window.ScrollRow = 13
$window.ScrollColumn = 1
$ws.Range("G20").Select()
